$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.4305099655428042
$ws.Range("D2").Value = 0.4289701383089291
$ws.Range("E2").Value = 0.4275098360405326

$ws.Range("C3").Value = 0.08444725529611972
$ws.Range("D3").Value = 0.08885734646265671
$ws.Range("E3").Value = 0.08941429023735328

$ws.Range("C4").Value = 0.007131338927048022
$ws.Range("D4").Value = 0.007895628020384611
$ws.Range("E4").Value = 0.00799491529864965

$ws.Range("C5").Value = 19.61563309914214
$ws.Range("D5").Value = 20.71411003407999
$ws.Range("E5").Value = 20.91514222584479
